$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper approach: for numeric-looking strings that must remain text (matching the
# source inlineStr cells), temporarily force a Text number format, assign the value,
# then restore the original style so no visible/style diff is introduced.

$ws.Range("D2").Value = "57.144.11"
$ws.Range("E2").Value = "  +1.58%  "

$ws.Range("D3").Value = "3.256.95"
$ws.Range("E3").Value = "  +0.89%  "

$ws.Range("E4").Value = "  +0.08%  "

$origStyle = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "397.53"
$ws.Range("D5").Style = $origStyle
$ws.Range("E5").Value = "  -0.48%  "

$origStyle = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "108.86"
$ws.Range("D6").Style = $origStyle
$ws.Range("E6").Value = "  -1.67%  "

$ws.Range("E7").Value = "  +4.02%  "

$origStyle = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.619"
$ws.Range("D9").Style = $origStyle
$ws.Range("E9").Value = "  -1.06%  "

$origStyle = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.25"
$ws.Range("D10").Style = $origStyle
$ws.Range("E10").Value = "  -0.67%  "

$origStyle = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0953"
$ws.Range("D11").Style = $origStyle

$ws.Range("E12").Value = "  +1.60%  "

$ws.Range("D13").Value = "3.773.25"
$ws.Range("E13").Value = "  +0.95%  "

$ws.Range("E14").Value = "  +0.93%  "

$origStyle = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "18.94"
$ws.Range("D15").Style = $origStyle
$ws.Range("E15").Value = "  -0.85%  "

$ws.Range("D16").Value = "3.252.40"
$ws.Range("E16").Value = "  +0.80%  "

$ws.Range("E17").Value = "  -2.35%  "

$origStyle = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "10.99"
$ws.Range("D18").Style = $origStyle
$ws.Range("E18").Value = "  +2.62%  "

$ws.Range("D19").Value = "56.949.98"
$ws.Range("E19").Value = "  +1.58%  "

$ws.Range("E20").Value = "  -1.66%  "

$ws.Range("E21").Value = "  +4.85%  "

$origStyle = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "12.92"
$ws.Range("D22").Style = $origStyle
$ws.Range("E22").Value = "  -1.33%  "

$origStyle = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "294.89"
$ws.Range("D23").Style = $origStyle
$ws.Range("E23").Value = "  -2.96%  "

$origStyle = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "73.92"
$ws.Range("D24").Style = $origStyle
$ws.Range("E24").Value = "  -2.06%  "

$origStyle = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.18"
$ws.Range("D25").Style = $origStyle
$ws.Range("E25").Value = "  -2.01%  "

$origStyle = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "28.06"
$ws.Range("D26").Style = $origStyle
$ws.Range("E26").Value = "  -0.82%  "

$ws.Range("B27").Value = "Filecoin"
$ws.Range("C27").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$origStyle = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.89"
$ws.Range("D27").Style = $origStyle
$ws.Range("E27").Value = "  -4.06%  "

$ws.Range("B28").Value = "LEO"
$ws.Range("C28").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$origStyle = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "4.38"
$ws.Range("D28").Style = $origStyle
$ws.Range("E28").Value = "  +0.38%  "

$origStyle = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.42"
$ws.Range("D29").Style = $origStyle
$ws.Range("E29").Value = "  -0.92%  "

$origStyle = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.169"
$ws.Range("D30").Style = $origStyle
$ws.Range("E30").Value = "  -3.00%  "

$ws.Range("E31").Value = "  +0.03%  "

$origStyle = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.113"
$ws.Range("D32").Style = $origStyle
$ws.Range("E32").Value = "  +1.15%  "

$origStyle = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "11.18"
$ws.Range("D33").Style = $origStyle
$ws.Range("E33").Value = "  -0.49%  "

$origStyle = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "40.68"
$ws.Range("D34").Style = $origStyle
$ws.Range("E34").Value = "  +12.13%  "

$origStyle = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0490"
$ws.Range("D35").Style = $origStyle
$ws.Range("E35").Value = "  -0.52%  "

$ws.Range("E36").Value = "  +1.18%  "

$origStyle = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "51.27"
$ws.Range("D37").Style = $origStyle
$ws.Range("E37").Value = "  -0.29%  "

$origStyle = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.999"
$ws.Range("D38").Style = $origStyle
$ws.Range("E38").Value = "  +0.04%  "

$ws.Range("E39").Value = "  -1.44%  "

$origStyle = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.99"
$ws.Range("D40").Style = $origStyle
$ws.Range("E40").Value = "  -3.87%  "

$ws.Range("E41").Value = "  +0.67%  "

$ws.Range("E42").Value = "  +0.75%  "

$ws.Range("B43").Value = "ARBITRUM"
$ws.Range("C43").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$origStyle = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.87"
$ws.Range("D43").Style = $origStyle
$ws.Range("E43").Value = "  -3.21%  "

$ws.Range("B44").Value = "TheGraph"
$ws.Range("C44").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$origStyle = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.283"
$ws.Range("D44").Style = $origStyle
$ws.Range("E44").Value = "  -1.25%  "

$origStyle = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.90"
$ws.Range("D45").Style = $origStyle
$ws.Range("E45").Value = "  -3.36%  "

$origStyle = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "16.71"
$ws.Range("D46").Style = $origStyle
$ws.Range("E46").Value = "  -3.31%  "

$origStyle = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "22.33"
$ws.Range("D47").Style = $origStyle
$ws.Range("E47").Value = "  +0.10%  "

$ws.Range("E48").Value = "  +4.27%  "

$origStyle = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.47"
$ws.Range("D49").Style = $origStyle
$ws.Range("E49").Value = "  -0.68%  "

$ws.Range("D50").Value = "2.143.72"
$ws.Range("E50").Value = "  +0.01%  "

$origStyle = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.97"
$ws.Range("D51").Style = $origStyle
$ws.Range("E51").Value = "  -7.93%  "
